$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2062
$ws.Range("I19").Value = 4247.3076
$ws.Range("J19").Value = 641.55
$ws.Range("K19").Value = 4247.3076
$ws.Range("L19").Value = 641.55
$ws.Range("M19").Value = -4072.3076
$ws.Range("N19").Value = -991.55
$ws.Range("H94").Value = 6100.625
$ws.Range("I94").Value = 5143.5713
$ws.Range("K94").Value = 5143.5713
$ws.Range("M94").Value = -4692.5713
$ws.Range("H132").Value = 2697.8064
$ws.Range("I132").Value = 1986.6428
$ws.Range("J132").Value = 9335.333000000001
$ws.Range("K132").Value = 5959.928400000001
$ws.Range("L132").Value = 28005.999
$ws.Range("M132").Value = -3429.928400000001
$ws.Range("N132").Value = -33065.999
$ws.Range("H137").Value = 3643.7917
$ws.Range("I137").Value = 4478
$ws.Range("J137").Value = 1975.375
$ws.Range("K137").Value = 13434
$ws.Range("L137").Value = 5926.125
$ws.Range("M137").Value = -10884
$ws.Range("N137").Value = -11026.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1457.6086
$ws.Range("I2").Value = 1505.9524
$ws.Range("J2").Value = 950
$ws.Range("K2").Value = 1505.9524
$ws.Range("L2").Value = 950
$ws.Range("M2").Value = -1392.9524
$ws.Range("N2").Value = -1176
$ws.Range("H6").Value = 500250
$ws.Range("I6").Value = 500250
$ws.Range("K6").Value = 500250
$ws.Range("M6").Value = -500077
$ws.Range("H37").Value = 22249.5
$ws.Range("I37").Value = 19499
$ws.Range("K37").Value = 19499
$ws.Range("M37").Value = -19226
$ws.Range("H61").Value = 467940.4
$ws.Range("I61").Value = 359054.8
$ws.Range("J61").Value = 671193.4399999999
$ws.Range("K61").Value = 359054.8
$ws.Range("L61").Value = 671193.4399999999
$ws.Range("M61").Value = -358842.8
$ws.Range("N61").Value = -671617.4399999999
$ws.Range("H102").Value = 5634.4443
$ws.Range("I102").Value = 1900
$ws.Range("J102").Value = 6701.4287
$ws.Range("K102").Value = 1900
$ws.Range("L102").Value = 6701.4287
$ws.Range("M102").Value = -278
$ws.Range("N102").Value = -9945.4287
$ws.Range("H116").Value = 1457.6086
$ws.Range("I116").Value = 1505.9524
$ws.Range("J116").Value = 950
$ws.Range("K116").Value = 1505.9524
$ws.Range("L116").Value = 950
$ws.Range("M116").Value = 788.0476000000001
$ws.Range("N116").Value = -5538
$ws.Range("H136").Value = 467940.4
$ws.Range("I136").Value = 359054.8
$ws.Range("J136").Value = 671193.4399999999
$ws.Range("K136").Value = 1077164.4
$ws.Range("L136").Value = 2013580.32
$ws.Range("M136").Value = -1074614.4
$ws.Range("N136").Value = -2018680.32

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1457.6086
$ws.Range("I3").Value = 1505.9524
$ws.Range("J3").Value = 950
$ws.Range("K3").Value = 1505.9524
$ws.Range("L3").Value = 950
$ws.Range("M3").Value = -1391.9524
$ws.Range("N3").Value = -1178
$ws.Range("H99").Value = 5500.909
$ws.Range("I99").Value = 7844.2856
$ws.Range("J99").Value = 1400
$ws.Range("K99").Value = 7844.2856
$ws.Range("L99").Value = 1400
$ws.Range("M99").Value = -6346.2856
$ws.Range("N99").Value = -4396

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 893.5333000000001
$ws.Range("I16").Value = 837.8570999999999
$ws.Range("K16").Value = 837.8570999999999
$ws.Range("M16").Value = -550.8570999999999
$ws.Range("H31").Value = 2917.3022
$ws.Range("I31").Value = 1991.6129
$ws.Range("J31").Value = 5308.6665
$ws.Range("K31").Value = 1991.6129
$ws.Range("L31").Value = 5308.6665
$ws.Range("M31").Value = -1696.6129
$ws.Range("N31").Value = -5898.6665
$ws.Range("H34").Value = 2917.3022
$ws.Range("I34").Value = 1991.6129
$ws.Range("J34").Value = 5308.6665
$ws.Range("K34").Value = 1991.6129
$ws.Range("L34").Value = 5308.6665
$ws.Range("M34").Value = -1789.6129
$ws.Range("N34").Value = -5712.6665
$ws.Range("H113").Value = 893.5333000000001
$ws.Range("I113").Value = 837.8570999999999
$ws.Range("K113").Value = 837.8570999999999
$ws.Range("M113").Value = 1332.1429
$ws.Range("H134").Value = 1795
$ws.Range("I134").Value = 1098.8334
$ws.Range("J134").Value = 2839.25
$ws.Range("K134").Value = 3296.5002
$ws.Range("L134").Value = 8517.75
$ws.Range("M134").Value = -761.5001999999999
$ws.Range("N134").Value = -13587.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 994.8444
$ws.Range("I5").Value = 406.12
$ws.Range("J5").Value = 1730.75
$ws.Range("K5").Value = 1218.36
$ws.Range("L5").Value = 5192.25
$ws.Range("M5").Value = -1106.36
$ws.Range("N5").Value = -5416.25
$ws.Range("H113").Value = 625.96155
$ws.Range("I113").Value = 614.2857
$ws.Range("J113").Value = 675
$ws.Range("K113").Value = 1842.8571
$ws.Range("L113").Value = 2025
$ws.Range("M113").Value = 327.1428999999998
$ws.Range("N113").Value = -6365
$ws.Range("H116").Value = 4266.3687
$ws.Range("I116").Value = 293
$ws.Range("J116").Value = 5685.4287
$ws.Range("K116").Value = 879
$ws.Range("L116").Value = 17056.2861
$ws.Range("M116").Value = 2563
$ws.Range("N116").Value = -23940.2861
$ws.Range("H135").Value = 994.8444
$ws.Range("I135").Value = 406.12
$ws.Range("J135").Value = 1730.75
$ws.Range("K135").Value = 3655.08
$ws.Range("L135").Value = 15576.75
$ws.Range("M135").Value = -1120.08
$ws.Range("N135").Value = -20646.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1283.7142
$ws.Range("I122").Value = 1213.8667
$ws.Range("J122").Value = 1458.3334
$ws.Range("K122").Value = 3641.6001
$ws.Range("L122").Value = 4375.0002
$ws.Range("M122").Value = -1191.6001
$ws.Range("N122").Value = -9275.0002
$ws.Range("H126").Value = 1829.95
$ws.Range("I126").Value = 1751.6296
$ws.Range("J126").Value = 1992.6154
$ws.Range("K126").Value = 5254.8888
$ws.Range("L126").Value = 5977.8462
$ws.Range("M126").Value = -2784.8888
$ws.Range("N126").Value = -10917.8462

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").ClearContents()
$ws.Range("N9").Value = 0
$ws.Range("H136").Value = 14869501
$ws.Range("I136").Value = 21762284
$ws.Range("J136").Value = 457319.53
$ws.Range("K136").Value = 65286852
$ws.Range("L136").Value = 1371958.59
$ws.Range("M136").Value = -65284302
$ws.Range("N136").Value = -1377058.59

Write-Host "Applied all changes"